# Extract-Skills stop-word cleanup: CV keyword cluster counts were
# recomputed after eliminating stop words (ip, it, ap, routers, video, ...).
# This updates the CLUSTER/SCORE fields in row 2 and the CV KEYWORDS list
# in column E, then drops the now-unused trailing keyword rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: SCORE / CLUSTER MATCH / CLUSTER SCORE fields ---
# Leading "'" keeps numeric-looking text (e.g. "8.333") stored as text,
# matching the workbook convention where these are shared strings, not numbers.
$ws.Range('D2').Value = '''8.333'   # SCORE
$ws.Range('F2').Value = 'analysis : 1'   # CLUSTER MUST HAVE MATCH
$ws.Range('G2').Value = '''9.090'   # CLUSTER MUST HAVE SCORE
$ws.Range('H2').ClearContents()   # CLUSTER GOOD TO HAVE MATCH no longer present
$ws.Range('I2').Value = '''0.0'   # CLUSTER GOOD TO HAVE SCORE
$ws.Range('K2').Value = '''0.0'   # CLUSTER SOFT SCORE

# --- Column E: recomputed CV KEYWORDS list (rows 2-100) ---
$ws.Range('E2').Value = 'networking : 3'
$ws.Range('E3').Value = 'cloud : 2'
$ws.Range('E4').Value = 'infrastructure : 1'
$ws.Range('E5').Value = 'english : 1'
$ws.Range('E6').Value = 'cisco : 6'
$ws.Range('E7').Value = 'cisco ccna : 4'
$ws.Range('E8').Value = 'ccna : 4'
$ws.Range('E9').Value = 'routing : 2'
$ws.Range('E10').Value = 'switching : 2'
$ws.Range('E11').Value = 'ccnp collaboration : 1'
$ws.Range('E12').Value = 'collaboration : 5'
$ws.Range('E13').Value = 'administrator : 3'
$ws.Range('E14').Value = 'network : 8'
$ws.Range('E15').Value = 'wifi : 1'
$ws.Range('E16').Value = 'monitoring : 4'
$ws.Range('E17').Value = 'process : 1'
$ws.Range('E18').Value = 'protocols : 1'
$ws.Range('E19').Value = 'management : 4'
$ws.Range('E20').Value = 'documentation : 1'
$ws.Range('E21').Value = 'risk assessment : 1'
$ws.Range('E22').Value = 'network administrator : 1'
$ws.Range('E23').Value = 'improvement : 1'
$ws.Range('E24').Value = 'administration : 2'
$ws.Range('E25').Value = 'security : 6'
$ws.Range('E26').Value = 'disaster recovery : 1'
$ws.Range('E27').Value = 'recovery : 1'
$ws.Range('E28').Value = 'trading : 1'
$ws.Range('E29').Value = 'access : 1'
$ws.Range('E30').Value = 'systems : 5'
$ws.Range('E31').Value = 'hybrid : 1'
$ws.Range('E32').Value = 'servers : 6'
$ws.Range('E33').Value = 'azure : 5'
$ws.Range('E34').Value = 'technical support : 6'
$ws.Range('E35').Value = 'migration : 6'
$ws.Range('E36').Value = 'unified communications : 1'
$ws.Range('E37').Value = 'communications : 1'
$ws.Range('E38').Value = 'cluster : 1'
$ws.Range('E39').Value = 'integration : 2'
$ws.Range('E40').Value = 'mobility : 1'
$ws.Range('E41').Value = 'tuning : 2'
$ws.Range('E42').Value = 'data center : 1'
$ws.Range('E43').Value = 'switches : 5'
$ws.Range('E44').Value = 'analog : 3'
$ws.Range('E45').Value = 'analysis : 2'
$ws.Range('E46').Value = 'troubleshooting : 1'
$ws.Range('E47').Value = 'engineer : 1'
$ws.Range('E48').Value = 'director : 1'
$ws.Range('E49').Value = 'global : 1'
$ws.Range('E50').Value = 'planning : 1'
$ws.Range('E51').Value = 'design : 3'
$ws.Range('E52').Value = 'hardware : 2'
$ws.Range('E53').Value = 'software support : 2'
$ws.Range('E54').Value = 'research : 1'
$ws.Range('E55').Value = 'vlan : 3'
$ws.Range('E56').Value = 'os : 2'
$ws.Range('E57').Value = 'virtualization : 2'
$ws.Range('E58').Value = 'backup : 1'
$ws.Range('E59').Value = 'sql : 1'
$ws.Range('E60').Value = 'microsoft server : 2'
$ws.Range('E61').Value = 'installation : 1'
$ws.Range('E62').Value = 'manager : 1'
$ws.Range('E63').Value = 'snmp : 2'
$ws.Range('E64').Value = 'zabbix : 1'
$ws.Range('E65').Value = 'uccx : 2'
$ws.Range('E66').Value = 'firewalls : 1'
$ws.Range('E67').Value = 'access control : 1'
$ws.Range('E68').Value = 'unity : 2'
$ws.Range('E69').Value = 'cube : 2'
$ws.Range('E70').Value = 'core : 1'
$ws.Range('E71').Value = 'acl : 2'
$ws.Range('E72').Value = 'nat : 2'
$ws.Range('E73').Value = 'ipsec : 2'
$ws.Range('E74').Value = 'ivr : 1'
$ws.Range('E75').Value = 'azure cloud : 1'
$ws.Range('E76').Value = 'firewall : 1'
$ws.Range('E77').Value = 'eigrp : 1'
$ws.Range('E78').Value = 'ospf : 1'
$ws.Range('E79').Value = 'vpn : 1'
$ws.Range('E80').Value = 'frame relay : 1'
$ws.Range('E81').Value = 'hsrp : 1'
$ws.Range('E82').Value = 'ssl : 1'
$ws.Range('E83').Value = 'vtp : 1'
$ws.Range('E84').Value = 'stp : 1'
$ws.Range('E85').Value = 'ssh : 1'
$ws.Range('E86').Value = 'syslog : 1'
$ws.Range('E87').Value = 'ipv6 : 1'
$ws.Range('E88').Value = 'tcp : 1'
$ws.Range('E89').Value = 'udp : 1'
$ws.Range('E90').Value = 'ftp : 1'
$ws.Range('E91').Value = 'smtp : 1'
$ws.Range('E92').Value = 'http : 1'
$ws.Range('E93').Value = 'https : 1'
$ws.Range('E94').Value = 'sip : 1'
$ws.Range('E95').Value = 'dhcp : 1'
$ws.Range('E96').Value = 'dns : 1'
$ws.Range('E97').Value = 'spark : 1'
$ws.Range('E98').Value = 'ad : 1'
$ws.Range('E99').Value = 'sql server : 1'
$ws.Range('E100').Value = 'it manager : 1'

# --- Drop the now-obsolete trailing keyword rows 101-115 ---
# (new keyword list only fills through row 100; dimension becomes A1:K100)
$ws.Range('A101:K115').EntireRow.Delete()
